# New weekly price record for "Feria Lagunitas de Puerto Montt - Espinaca".
# Insert a new row above the existing row 5, shifting rows 5:27 down to
# 6:28 (dimension grows from A1:R27 to A1:R28), then populate the new
# row 5 with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a whole new row at position 5.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the latest weekly observation.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44750
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112012
$ws.Range("G5").Value = "Espinaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("N5").Value = "`$/cuna 10 kilos"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 1300
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = "Hortaliza"
